$d = $word.ActiveDocument

$d.Content.Find.Execute('2025-06-12 Thursday', $true, $true, $false, $false, $false, $true, 1, $false, '2025-06-13 Friday', 2) | Out-Null
$d.Content.Find.Execute('9+66=', $true, $true, $false, $false, $false, $true, 1, $false, '80-47=', 2) | Out-Null
$d.Content.Find.Execute('8+11=', $true, $true, $false, $false, $false, $true, 1, $false, '10+52=', 2) | Out-Null
$d.Content.Find.Execute('40-30=', $true, $true, $false, $false, $false, $true, 1, $false, '44+43=', 2) | Out-Null
$d.Content.Find.Execute('1+10=', $true, $true, $false, $false, $false, $true, 1, $false, '84-7=', 2) | Out-Null
$d.Content.Find.Execute('60-14=', $true, $true, $false, $false, $false, $true, 1, $false, '9+2=', 2) | Out-Null
$d.Content.Find.Execute('25+11=', $true, $true, $false, $false, $false, $true, 1, $false, '15+8=', 2) | Out-Null
$d.Content.Find.Execute('55+17=', $true, $true, $false, $false, $false, $true, 1, $false, '77-23=', 2) | Out-Null
$d.Content.Find.Execute('73-34=', $true, $true, $false, $false, $false, $true, 1, $false, '41+43=', 2) | Out-Null
$d.Content.Find.Execute('84-1=', $true, $true, $false, $false, $false, $true, 1, $false, '71-54=', 2) | Out-Null
$d.Content.Find.Execute('98-90=', $true, $true, $false, $false, $false, $true, 1, $false, '42+7=', 2) | Out-Null
$d.Content.Find.Execute('11+79=', $true, $true, $false, $false, $false, $true, 1, $false, '91-55=', 2) | Out-Null
$d.Content.Find.Execute('43+35=', $true, $true, $false, $false, $false, $true, 1, $false, '38+24=', 2) | Out-Null
$d.Content.Find.Execute('36-7=', $true, $true, $false, $false, $false, $true, 1, $false, '72-25=', 2) | Out-Null
$d.Content.Find.Execute('11+8=', $true, $true, $false, $false, $false, $true, 1, $false, '5+22=', 2) | Out-Null
$d.Content.Find.Execute('19+5=', $true, $true, $false, $false, $false, $true, 1, $false, '95-91=', 2) | Out-Null
$d.Content.Find.Execute('82-9=', $true, $true, $false, $false, $false, $true, 1, $false, '73-59=', 2) | Out-Null
$d.Content.Find.Execute('14+39=', $true, $true, $false, $false, $false, $true, 1, $false, '80-0=', 2) | Out-Null
$d.Content.Find.Execute('52+9=', $true, $true, $false, $false, $false, $true, 1, $false, '93+1=', 2) | Out-Null
$d.Content.Find.Execute('91-50=', $true, $true, $false, $false, $false, $true, 1, $false, '9+37=', 2) | Out-Null
$d.Content.Find.Execute('37+21=', $true, $true, $false, $false, $false, $true, 1, $false, '44+37=', 2) | Out-Null
$d.Content.Find.Execute('98-66=', $true, $true, $false, $false, $false, $true, 1, $false, '10+16=', 2) | Out-Null
$d.Content.Find.Execute('24-16=', $true, $true, $false, $false, $false, $true, 1, $false, '54-7=', 2) | Out-Null
$d.Content.Find.Execute('55-19=', $true, $true, $false, $false, $false, $true, 1, $false, '58-14=', 2) | Out-Null
$d.Content.Find.Execute('71-32=', $true, $true, $false, $false, $false, $true, 1, $false, '77-34=', 2) | Out-Null
$d.Content.Find.Execute('49-24=', $true, $true, $false, $false, $false, $true, 1, $false, '22+48=', 2) | Out-Null
$d.Content.Find.Execute('83-9=', $true, $true, $false, $false, $false, $true, 1, $false, '67-8=', 2) | Out-Null
$d.Content.Find.Execute('84-78=', $true, $true, $false, $false, $false, $true, 1, $false, '72-1=', 2) | Out-Null
$d.Content.Find.Execute('65-52=', $true, $true, $false, $false, $false, $true, 1, $false, '41+54=', 2) | Out-Null
$d.Content.Find.Execute('20+58=', $true, $true, $false, $false, $false, $true, 1, $false, '27+15=', 2) | Out-Null
$d.Content.Find.Execute('34+24=', $true, $true, $false, $false, $false, $true, 1, $false, '72+16=', 2) | Out-Null
$d.Content.Find.Execute('68-39=', $true, $true, $false, $false, $false, $true, 1, $false, '26+4=', 2) | Out-Null
$d.Content.Find.Execute('16-14=', $true, $true, $false, $false, $false, $true, 1, $false, '35+34=', 2) | Out-Null
$d.Content.Find.Execute('37+35=', $true, $true, $false, $false, $false, $true, 1, $false, '21+32=', 2) | Out-Null
$d.Content.Find.Execute('46-26=', $true, $true, $false, $false, $false, $true, 1, $false, '89-67=', 2) | Out-Null
$d.Content.Find.Execute('80-35=', $true, $true, $false, $false, $false, $true, 1, $false, '45+12=', 2) | Out-Null
$d.Content.Find.Execute('5+32=', $true, $true, $false, $false, $false, $true, 1, $false, '66+27=', 2) | Out-Null
$d.Content.Find.Execute('40+5=', $true, $true, $false, $false, $false, $true, 1, $false, '75+23=', 2) | Out-Null
$d.Content.Find.Execute('15+24=', $true, $true, $false, $false, $false, $true, 1, $false, '23+53=', 2) | Out-Null
$d.Content.Find.Execute('61-12=', $true, $true, $false, $false, $false, $true, 1, $false, '79+6=', 2) | Out-Null
$d.Content.Find.Execute('78-3=', $true, $true, $false, $false, $false, $true, 1, $false, '15+79=', 2) | Out-Null
$d.Content.Find.Execute('15+34=', $true, $true, $false, $false, $false, $true, 1, $false, '74-62=', 2) | Out-Null
$d.Content.Find.Execute('96+3=', $true, $true, $false, $false, $false, $true, 1, $false, '12+61=', 2) | Out-Null
$d.Content.Find.Execute('89-44=', $true, $true, $false, $false, $false, $true, 1, $false, '13+66=', 2) | Out-Null
$d.Content.Find.Execute('27+29=', $true, $true, $false, $false, $false, $true, 1, $false, '35+37=', 2) | Out-Null
$d.Content.Find.Execute('53+40=', $true, $true, $false, $false, $false, $true, 1, $false, '51-33=', 2) | Out-Null
$d.Content.Find.Execute('40+38=', $true, $true, $false, $false, $false, $true, 1, $false, '43+54=', 2) | Out-Null
$d.Content.Find.Execute('45+4=', $true, $true, $false, $false, $false, $true, 1, $false, '1+96=', 2) | Out-Null
$d.Content.Find.Execute('58-45=', $true, $true, $false, $false, $false, $true, 1, $false, '23-17=', 2) | Out-Null
$d.Content.Find.Execute('39+33=', $true, $true, $false, $false, $false, $true, 1, $false, '25+58=', 2) | Out-Null
$d.Content.Find.Execute('7+84=', $true, $true, $false, $false, $false, $true, 1, $false, '80+8=', 2) | Out-Null
$d.Content.Find.Execute('19-13=', $true, $true, $false, $false, $false, $true, 1, $false, '10-7=', 2) | Out-Null
$d.Content.Find.Execute('50+38=', $true, $true, $false, $false, $false, $true, 1, $false, '88-67=', 2) | Out-Null
$d.Content.Find.Execute('83-13=', $true, $true, $false, $false, $false, $true, 1, $false, '31+45=', 2) | Out-Null
$d.Content.Find.Execute('48+34=', $true, $true, $false, $false, $false, $true, 1, $false, '52+33=', 2) | Out-Null
$d.Content.Find.Execute('2+34=', $true, $true, $false, $false, $false, $true, 1, $false, '37-23=', 2) | Out-Null
$d.Content.Find.Execute('75+7=', $true, $true, $false, $false, $false, $true, 1, $false, '92-78=', 2) | Out-Null
$d.Content.Find.Execute('31-27=', $true, $true, $false, $false, $false, $true, 1, $false, '41+54=', 2) | Out-Null
$d.Content.Find.Execute('95-38=', $true, $true, $false, $false, $false, $true, 1, $false, '0+39=', 2) | Out-Null
$d.Content.Find.Execute('44-25=', $true, $true, $false, $false, $false, $true, 1, $false, '33+23=', 2) | Out-Null
$d.Content.Find.Execute('53-38=', $true, $true, $false, $false, $false, $true, 1, $false, '41-3=', 2) | Out-Null
$d.Content.Find.Execute('42+34=', $true, $true, $false, $false, $false, $true, 1, $false, '73-18=', 2) | Out-Null
$d.Content.Find.Execute('84-57=', $true, $true, $false, $false, $false, $true, 1, $false, '57-38=', 2) | Out-Null
$d.Content.Find.Execute('21+5=', $true, $true, $false, $false, $false, $true, 1, $false, '83-73=', 2) | Out-Null
$d.Content.Find.Execute('5+89=', $true, $true, $false, $false, $false, $true, 1, $false, '29+67=', 2) | Out-Null
$d.Content.Find.Execute('94+0=', $true, $true, $false, $false, $false, $true, 1, $false, '50-38=', 2) | Out-Null
$d.Content.Find.Execute('45+1=', $true, $true, $false, $false, $false, $true, 1, $false, '37+42=', 2) | Out-Null
$d.Content.Find.Execute('11+45=', $true, $true, $false, $false, $false, $true, 1, $false, '67-39=', 2) | Out-Null
$d.Content.Find.Execute('66+16=', $true, $true, $false, $false, $false, $true, 1, $false, '49+7=', 2) | Out-Null
$d.Content.Find.Execute('51+42=', $true, $true, $false, $false, $false, $true, 1, $false, '75-47=', 2) | Out-Null
$d.Content.Find.Execute('51-47=', $true, $true, $false, $false, $false, $true, 1, $false, '50-32=', 2) | Out-Null
$d.Content.Find.Execute('24+10=', $true, $true, $false, $false, $false, $true, 1, $false, '32-24=', 2) | Out-Null
$d.Content.Find.Execute('73-51=', $true, $true, $false, $false, $false, $true, 1, $false, '49-49=', 2) | Out-Null
$d.Content.Find.Execute('43+46=', $true, $true, $false, $false, $false, $true, 1, $false, '21+26=', 2) | Out-Null
$d.Content.Find.Execute('18+26=', $true, $true, $false, $false, $false, $true, 1, $false, '37-3=', 2) | Out-Null
$d.Content.Find.Execute('49-29=', $true, $true, $false, $false, $false, $true, 1, $false, '8-6=', 2) | Out-Null
$d.Content.Find.Execute('73+16=', $true, $true, $false, $false, $false, $true, 1, $false, '92-20=', 2) | Out-Null
$d.Content.Find.Execute('46-44=', $true, $true, $false, $false, $false, $true, 1, $false, '8+75=', 2) | Out-Null
$d.Content.Find.Execute('77-73=', $true, $true, $false, $false, $false, $true, 1, $false, '20+8=', 2) | Out-Null
$d.Content.Find.Execute('23+35=', $true, $true, $false, $false, $false, $true, 1, $false, '18+60=', 2) | Out-Null
$d.Content.Find.Execute('24+64=', $true, $true, $false, $false, $false, $true, 1, $false, '89-21=', 2) | Out-Null
$d.Content.Find.Execute('76-43=', $true, $true, $false, $false, $false, $true, 1, $false, '63-27=', 2) | Out-Null
$d.Content.Find.Execute('40-24=', $true, $true, $false, $false, $false, $true, 1, $false, '11+59=', 2) | Out-Null
$d.Content.Find.Execute('76-30=', $true, $true, $false, $false, $false, $true, 1, $false, '18+41=', 2) | Out-Null
$d.Content.Find.Execute('90-16=', $true, $true, $false, $false, $false, $true, 1, $false, '36+52=', 2) | Out-Null
$d.Content.Find.Execute('71-68=', $true, $true, $false, $false, $false, $true, 1, $false, '56+6=', 2) | Out-Null
$d.Content.Find.Execute('58-10=', $true, $true, $false, $false, $false, $true, 1, $false, '73-24=', 2) | Out-Null
$d.Content.Find.Execute('58-12=', $true, $true, $false, $false, $false, $true, 1, $false, '59+40=', 2) | Out-Null
$d.Content.Find.Execute('60+33=', $true, $true, $false, $false, $false, $true, 1, $false, '2+2=', 2) | Out-Null
$d.Content.Find.Execute('27+49=', $true, $true, $false, $false, $false, $true, 1, $false, '85+10=', 2) | Out-Null
$d.Content.Find.Execute('57-42=', $true, $true, $false, $false, $false, $true, 1, $false, '84-40=', 2) | Out-Null
$d.Content.Find.Execute('12-5=', $true, $true, $false, $false, $false, $true, 1, $false, '93-15=', 2) | Out-Null
$d.Content.Find.Execute('71+18=', $true, $true, $false, $false, $false, $true, 1, $false, '53+35=', 2) | Out-Null
$d.Content.Find.Execute('10+12=', $true, $true, $false, $false, $false, $true, 1, $false, '72-18=', 2) | Out-Null
$d.Content.Find.Execute('84+13=', $true, $true, $false, $false, $false, $true, 1, $false, '60-22=', 2) | Out-Null
$d.Content.Find.Execute('86-6=', $true, $true, $false, $false, $false, $true, 1, $false, '26+11=', 2) | Out-Null
$d.Content.Find.Execute('86-85=', $true, $true, $false, $false, $false, $true, 1, $false, '1+31=', 2) | Out-Null
$d.Content.Find.Execute('0+30=', $true, $true, $false, $false, $false, $true, 1, $false, '14-2=', 2) | Out-Null
$d.Content.Find.Execute('51-40=', $true, $true, $false, $false, $false, $true, 1, $false, '35-11=', 2) | Out-Null
$d.Content.Find.Execute('75-23=', $true, $true, $false, $false, $false, $true, 1, $false, '2+9=', 2) | Out-Null
$d.Content.Find.Execute('85-34=', $true, $true, $false, $false, $false, $true, 1, $false, '36-16=', 2) | Out-Null
